$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 265.83334
$ws.Range("I9").Value = 209.1
$ws.Range("K9").Value = 209.1
$ws.Range("M9").Value = -40.09999999999999
$ws.Range("H12").Value = 485.5
$ws.Range("J12").Value = 880
$ws.Range("L12").Value = 880
$ws.Range("N12").Value = -1220
$ws.Range("H18").Value = 2636.8333
$ws.Range("I18").Value = 1972.75
$ws.Range("J18").Value = 3965
$ws.Range("K18").Value = 1972.75
$ws.Range("L18").Value = 3965
$ws.Range("M18").Value = -1688.75
$ws.Range("N18").Value = -4533
$ws.Range("H19").Value = 5883469.5
$ws.Range("J19").Value = 10001022
$ws.Range("L19").Value = 10001022
$ws.Range("N19").Value = -10001372
$ws.Range("H42").Value = 1329.75
$ws.Range("I42").Value = 459
$ws.Range("K42").Value = 1377
$ws.Range("M42").Value = -1147
$ws.Range("H112").Value = 1923.6786
$ws.Range("J112").Value = 1937.1923
$ws.Range("L112").Value = 5811.5769
$ws.Range("N112").Value = -8027.5769
$ws.Range("H135").Value = 4668.1816
$ws.Range("I135").Value = 907.1429000000001
$ws.Range("K135").Value = 8164.2861
$ws.Range("M135").Value = -5629.2861
$ws.Range("H137").Value = 25003568
$ws.Range("I137").Value = 31253654
$ws.Range("J137").Value = 3225
$ws.Range("K137").Value = 93760962
$ws.Range("L137").Value = 9675
$ws.Range("M137").Value = -93758412
$ws.Range("N137").Value = -14775

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 266.27274
$ws.Range("I5").Value = 288.66666
$ws.Range("J5").Value = 239.4
$ws.Range("K5").Value = 288.66666
$ws.Range("L5").Value = 239.4
$ws.Range("M5").Value = -176.66666
$ws.Range("N5").Value = -463.4
$ws.Range("H74").Value = 2254.3936
$ws.Range("I74").Value = 2312.26
$ws.Range("J74").Value = 1991.3636
$ws.Range("K74").Value = 2312.26
$ws.Range("L74").Value = 1991.3636
$ws.Range("M74").Value = -1438.26
$ws.Range("N74").Value = -3739.3636
$ws.Range("H77").Value = 2254.3936
$ws.Range("I77").Value = 2312.26
$ws.Range("J77").Value = 1991.3636
$ws.Range("K77").Value = 11561.3
$ws.Range("L77").Value = 9956.817999999999
$ws.Range("M77").Value = -7193.300000000001
$ws.Range("N77").Value = -18692.818

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 266.27274
$ws.Range("I4").Value = 288.66666
$ws.Range("J4").Value = 239.4
$ws.Range("K4").Value = 288.66666
$ws.Range("L4").Value = 239.4
$ws.Range("M4").Value = -173.66666
$ws.Range("N4").Value = -469.4
$ws.Range("H122").Value = 0
$ws.Range("J122").Value = 0
$ws.Range("L122").Value = 0
$ws.Range("N122").ClearContents()
$ws.Range("J127").Value = 0
$ws.Range("L127").Value = 0
$ws.Range("N127").ClearContents()
$ws.Range("H134").Value = 4350033
$ws.Range("I134").Value = 2226.9524
$ws.Range("K134").Value = 6680.8572
$ws.Range("M134").Value = -4145.8572

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 14.142858
$ws.Range("I7").Value = 14.5
$ws.Range("K7").Value = 14.5
$ws.Range("M7").Value = 98.5
$ws.Range("H35").Value = 13532.667
$ws.Range("I35").Value = 10299
$ws.Range("K35").Value = 10299
$ws.Range("M35").Value = -10005
$ws.Range("H74").Value = 99999
$ws.Range("J74").Value = 99999
$ws.Range("L74").Value = 99999
$ws.Range("N74").Value = -101747
$ws.Range("H77").Value = 99999
$ws.Range("J77").Value = 99999
$ws.Range("L77").Value = 299997
$ws.Range("N77").Value = -308733
$ws.Range("H97").Value = 59598
$ws.Range("J97").Value = 59598
$ws.Range("L97").Value = 59598
$ws.Range("N97").Value = -61580
$ws.Range("H99").Value = 14482.667
$ws.Range("I99").Value = 6026.9
$ws.Range("K99").Value = 6026.9
$ws.Range("M99").Value = -4528.9
$ws.Range("H105").Value = 1521.25
$ws.Range("J105").Value = 4500
$ws.Range("L105").Value = 4500
$ws.Range("N105").Value = -7994
$ws.Range("H126").Value = 14482.667
$ws.Range("I126").Value = 6026.9
$ws.Range("K126").Value = 18080.7
$ws.Range("M126").Value = -15610.7
$ws.Range("H132").Value = 3673
$ws.Range("I132").Value = 1762.091
$ws.Range("K132").Value = 5286.272999999999
$ws.Range("M132").Value = -2756.272999999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H33").Value = 5330834
$ws.Range("I33").Value = 145
$ws.Range("J33").Value = 7700029
$ws.Range("K33").Value = 870
$ws.Range("L33").Value = 46200174
$ws.Range("M33").Value = -587
$ws.Range("N33").Value = -46200740
$ws.Range("H120").Value = 31099.666
$ws.Range("I120").Value = 29967
$ws.Range("J120").Value = 31666
$ws.Range("K120").Value = 89901
$ws.Range("L120").Value = 94998
$ws.Range("M120").Value = -85063
$ws.Range("N120").Value = -104674
$ws.Range("H121").Value = 2580.7715
$ws.Range("J121").Value = 3118.2144
$ws.Range("L121").Value = 9354.643199999999
$ws.Range("N121").Value = -11974.6432
$ws.Range("H131").Value = 3621.0881
$ws.Range("I131").Value = 2478.2144
$ws.Range("K131").Value = 7434.6432
$ws.Range("M131").Value = -2394.6432
$ws.Range("H132").Value = 3796.2222
$ws.Range("I132").Value = 2474.25
$ws.Range("J132").Value = 4853.8
$ws.Range("K132").Value = 22268.25
$ws.Range("L132").Value = 43684.2
$ws.Range("M132").Value = -19738.25
$ws.Range("N132").Value = -48744.2
$ws.Range("H138").Value = 13483.583
$ws.Range("I138").Value = 12347
$ws.Range("K138").Value = 37041
$ws.Range("M138").Value = -31901

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 15324.875
$ws.Range("I22").Value = 38599.668
$ws.Range("K22").Value = 38599.668
$ws.Range("M22").Value = -38304.668
$ws.Range("H27").Value = 15324.875
$ws.Range("I27").Value = 38599.668
$ws.Range("K27").Value = 38599.668
$ws.Range("M27").Value = -38492.668
$ws.Range("H40").Value = 7041.5713
$ws.Range("I40").Value = 6464
$ws.Range("J40").Value = 7474.75
$ws.Range("K40").Value = 6464
$ws.Range("L40").Value = 7474.75
$ws.Range("M40").Value = -6328
$ws.Range("N40").Value = -7746.75
$ws.Range("H61").Value = 66671052
$ws.Range("I61").Value = 125003270
$ws.Range("K61").Value = 125003270
$ws.Range("M61").Value = -125003068
$ws.Range("H113").Value = 66671052
$ws.Range("I113").Value = 125003270
$ws.Range("K113").Value = 125003270
$ws.Range("M113").Value = -125001100
$ws.Range("H122").Value = 3895.3408
$ws.Range("I122").Value = 3499.9429
$ws.Range("J122").Value = 5433
$ws.Range("K122").Value = 10499.8287
$ws.Range("L122").Value = 16299
$ws.Range("M122").Value = -8049.8287
$ws.Range("N122").Value = -21199
$ws.Range("H132").Value = 3188.7144
$ws.Range("I132").Value = 1732.7391
$ws.Range("K132").Value = 5198.2173
$ws.Range("M132").Value = -2668.2173
$ws.Range("H136").Value = 1958.2821
$ws.Range("I136").Value = 1213.3448
$ws.Range("K136").Value = 3640.0344
$ws.Range("M136").Value = -1090.0344

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 71434904
$ws.Range("I96").Value = 5877.7
$ws.Range("K96").Value = 5877.7
$ws.Range("M96").Value = -4504.7
$ws.Range("H100").Value = 943.1429000000001
$ws.Range("I100").Value = 1039
$ws.Range("K100").Value = 2078
$ws.Range("M100").Value = -1537
$ws.Range("H104").Value = 118842.25
$ws.Range("J104").Value = 118842.25
$ws.Range("L104").Value = 118842.25
$ws.Range("N104").Value = -125830.25
$ws.Range("H106").Value = 99999
$ws.Range("J106").Value = 99999
$ws.Range("L106").Value = 99999
$ws.Range("N106").Value = -102523
$ws.Range("H132").Value = 836013.4399999999
$ws.Range("I132").Value = 2906.889
$ws.Range("K132").Value = 8720.667000000001
$ws.Range("M132").Value = -6190.667000000001
